$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.065099835395813
$ws.Range("B1").Value = 1.671149730682373
$ws.Range("C1").Value = 6.886199474334717
$ws.Range("D1").Value = 2.715194940567017
$ws.Range("E1").Value = 1.473120331764221
